$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh particelle_non_trovate data (rows 18-88): update codice_particella / codice_comune_catastale values
$ws.Cells.Item(18, 2).Value = '8974/1'
$ws.Cells.Item(18, 3).Value = 154
$ws.Cells.Item(19, 2).NumberFormat = "@"
$ws.Cells.Item(19, 2).Value = '15380'
$ws.Cells.Item(19, 2).Style = "Normal"
$ws.Cells.Item(19, 3).Value = 154
$ws.Cells.Item(20, 2).NumberFormat = "@"
$ws.Cells.Item(20, 2).Value = '391'
$ws.Cells.Item(20, 2).Style = "Normal"
$ws.Cells.Item(20, 3).Value = 174
$ws.Cells.Item(21, 2).NumberFormat = "@"
$ws.Cells.Item(21, 2).Value = '466'
$ws.Cells.Item(21, 2).Style = "Normal"
$ws.Cells.Item(21, 3).Value = 213
$ws.Cells.Item(22, 2).NumberFormat = "@"
$ws.Cells.Item(22, 2).Value = '.315'
$ws.Cells.Item(22, 2).Style = "Normal"
$ws.Cells.Item(22, 3).Value = 282
$ws.Cells.Item(23, 2).NumberFormat = "@"
$ws.Cells.Item(23, 2).Value = '403'
$ws.Cells.Item(23, 2).Style = "Normal"
$ws.Cells.Item(23, 3).Value = 97
$ws.Cells.Item(24, 2).NumberFormat = "@"
$ws.Cells.Item(24, 2).Value = '746'
$ws.Cells.Item(24, 2).Style = "Normal"
$ws.Cells.Item(24, 3).Value = 97
$ws.Cells.Item(25, 2).NumberFormat = "@"
$ws.Cells.Item(25, 2).Value = '749'
$ws.Cells.Item(25, 2).Style = "Normal"
$ws.Cells.Item(25, 3).Value = 97
$ws.Cells.Item(26, 2).NumberFormat = "@"
$ws.Cells.Item(26, 2).Value = '754'
$ws.Cells.Item(26, 2).Style = "Normal"
$ws.Cells.Item(26, 3).Value = 97
$ws.Cells.Item(27, 2).NumberFormat = "@"
$ws.Cells.Item(27, 2).Value = '755'
$ws.Cells.Item(27, 2).Style = "Normal"
$ws.Cells.Item(27, 3).Value = 97
$ws.Cells.Item(28, 2).NumberFormat = "@"
$ws.Cells.Item(28, 2).Value = '757'
$ws.Cells.Item(28, 2).Style = "Normal"
$ws.Cells.Item(28, 3).Value = 97
$ws.Cells.Item(29, 2).NumberFormat = "@"
$ws.Cells.Item(29, 2).Value = '758'
$ws.Cells.Item(29, 2).Style = "Normal"
$ws.Cells.Item(29, 3).Value = 97
$ws.Cells.Item(30, 2).Value = '825/63'
$ws.Cells.Item(30, 3).Value = 97
$ws.Cells.Item(31, 2).Value = '1900/4'
$ws.Cells.Item(31, 3).Value = 317
$ws.Cells.Item(32, 2).Value = '1272/3'
$ws.Cells.Item(32, 3).Value = 390
$ws.Cells.Item(33, 2).NumberFormat = "@"
$ws.Cells.Item(33, 2).Value = '765'
$ws.Cells.Item(33, 2).Style = "Normal"
$ws.Cells.Item(33, 3).Value = 404
$ws.Cells.Item(34, 2).Value = '94/6'
$ws.Cells.Item(34, 3).Value = 251
$ws.Cells.Item(35, 2).Value = '2681/1'
$ws.Cells.Item(35, 3).Value = 442
$ws.Cells.Item(36, 2).NumberFormat = "@"
$ws.Cells.Item(36, 2).Value = '789'
$ws.Cells.Item(36, 2).Style = "Normal"
$ws.Cells.Item(36, 3).Value = 443
$ws.Cells.Item(37, 2).NumberFormat = "@"
$ws.Cells.Item(37, 2).Value = '53'
$ws.Cells.Item(37, 2).Style = "Normal"
$ws.Cells.Item(37, 3).Value = 215
$ws.Cells.Item(38, 2).NumberFormat = "@"
$ws.Cells.Item(38, 2).Value = '454'
$ws.Cells.Item(38, 2).Style = "Normal"
$ws.Cells.Item(38, 3).Value = 215
$ws.Cells.Item(39, 2).Value = '420/80'
$ws.Cells.Item(39, 3).Value = 215
$ws.Cells.Item(40, 2).Value = '420/92'
$ws.Cells.Item(40, 3).Value = 215
$ws.Cells.Item(41, 2).Value = '420/93'
$ws.Cells.Item(41, 3).Value = 215
$ws.Cells.Item(42, 2).Value = '420/94'
$ws.Cells.Item(42, 3).Value = 215
$ws.Cells.Item(43, 2).Value = '420/95'
$ws.Cells.Item(43, 3).Value = 215
$ws.Cells.Item(44, 2).Value = '420/96'
$ws.Cells.Item(44, 3).Value = 215
$ws.Cells.Item(45, 2).Value = '420/97'
$ws.Cells.Item(45, 3).Value = 215
$ws.Cells.Item(46, 2).Value = '420/101'
$ws.Cells.Item(46, 3).Value = 215
$ws.Cells.Item(47, 2).Value = '420/102'
$ws.Cells.Item(47, 3).Value = 215
$ws.Cells.Item(48, 2).Value = '420/106'
$ws.Cells.Item(48, 3).Value = 215
$ws.Cells.Item(49, 2).Value = '420/107'
$ws.Cells.Item(49, 3).Value = 215
$ws.Cells.Item(50, 2).Value = '420/109'
$ws.Cells.Item(50, 3).Value = 215
$ws.Cells.Item(51, 2).Value = '420/110'
$ws.Cells.Item(51, 3).Value = 215
$ws.Cells.Item(52, 2).Value = '705/11'
$ws.Cells.Item(52, 3).Value = 215
$ws.Cells.Item(53, 2).NumberFormat = "@"
$ws.Cells.Item(53, 2).Value = '756'
$ws.Cells.Item(53, 2).Style = "Normal"
$ws.Cells.Item(53, 3).Value = 215
$ws.Cells.Item(54, 2).Value = '798/3'
$ws.Cells.Item(54, 3).Value = 215
$ws.Cells.Item(55, 2).Value = '1411/1'
$ws.Cells.Item(55, 3).Value = 256
$ws.Cells.Item(56, 2).Value = '1411/2'
$ws.Cells.Item(56, 3).Value = 256
$ws.Cells.Item(57, 2).Value = '1411/3'
$ws.Cells.Item(57, 3).Value = 256
$ws.Cells.Item(58, 2).Value = '1411/4'
$ws.Cells.Item(58, 3).Value = 256
$ws.Cells.Item(59, 2).Value = '1411/5'
$ws.Cells.Item(59, 3).Value = 256
$ws.Cells.Item(60, 2).NumberFormat = "@"
$ws.Cells.Item(60, 2).Value = '1412'
$ws.Cells.Item(60, 2).Style = "Normal"
$ws.Cells.Item(60, 3).Value = 256
$ws.Cells.Item(61, 2).NumberFormat = "@"
$ws.Cells.Item(61, 2).Value = '1488'
$ws.Cells.Item(61, 2).Style = "Normal"
$ws.Cells.Item(61, 3).Value = 256
$ws.Cells.Item(62, 2).Value = '1117/2'
$ws.Cells.Item(62, 3).Value = 193
$ws.Cells.Item(63, 2).Value = '1230/100'
$ws.Cells.Item(63, 3).Value = 193
$ws.Cells.Item(64, 2).Value = '1230/115'
$ws.Cells.Item(64, 3).Value = 193
$ws.Cells.Item(65, 2).Value = '1230/85'
$ws.Cells.Item(65, 3).Value = 193
$ws.Cells.Item(66, 2).Value = '1230/86'
$ws.Cells.Item(66, 3).Value = 193
$ws.Cells.Item(67, 2).Value = '1230/87'
$ws.Cells.Item(67, 3).Value = 193
$ws.Cells.Item(68, 2).Value = '1230/88'
$ws.Cells.Item(68, 3).Value = 193
$ws.Cells.Item(69, 2).Value = '1303/1'
$ws.Cells.Item(69, 3).Value = 193
$ws.Cells.Item(70, 2).Value = '1303/2'
$ws.Cells.Item(70, 3).Value = 193
$ws.Cells.Item(71, 2).NumberFormat = "@"
$ws.Cells.Item(71, 2).Value = '1309'
$ws.Cells.Item(71, 2).Style = "Normal"
$ws.Cells.Item(71, 3).Value = 193
$ws.Cells.Item(72, 2).NumberFormat = "@"
$ws.Cells.Item(72, 2).Value = '1330'
$ws.Cells.Item(72, 2).Style = "Normal"
$ws.Cells.Item(72, 3).Value = 193
$ws.Cells.Item(73, 2).NumberFormat = "@"
$ws.Cells.Item(73, 2).Value = '1334'
$ws.Cells.Item(73, 2).Style = "Normal"
$ws.Cells.Item(73, 3).Value = 193
$ws.Cells.Item(74, 2).NumberFormat = "@"
$ws.Cells.Item(74, 2).Value = '1346'
$ws.Cells.Item(74, 2).Style = "Normal"
$ws.Cells.Item(74, 3).Value = 193
$ws.Cells.Item(75, 2).Value = '1369/1'
$ws.Cells.Item(75, 3).Value = 193
$ws.Cells.Item(76, 2).Value = '194/4'
$ws.Cells.Item(76, 3).Value = 193
$ws.Cells.Item(77, 2).Value = '254/2'
$ws.Cells.Item(77, 3).Value = 193
$ws.Cells.Item(78, 2).Value = '337/5'
$ws.Cells.Item(78, 3).Value = 193
$ws.Cells.Item(79, 2).Value = '393/1'
$ws.Cells.Item(79, 3).Value = 193
$ws.Cells.Item(80, 2).Value = '393/2'
$ws.Cells.Item(80, 3).Value = 193
$ws.Cells.Item(81, 2).Value = '393/3'
$ws.Cells.Item(81, 3).Value = 193
$ws.Cells.Item(82, 2).NumberFormat = "@"
$ws.Cells.Item(82, 2).Value = '465'
$ws.Cells.Item(82, 2).Style = "Normal"
$ws.Cells.Item(82, 3).Value = 193
$ws.Cells.Item(83, 2).NumberFormat = "@"
$ws.Cells.Item(83, 2).Value = '614'
$ws.Cells.Item(83, 2).Style = "Normal"
$ws.Cells.Item(83, 3).Value = 193
$ws.Cells.Item(84, 2).Value = '384/1'
$ws.Cells.Item(84, 3).Value = 193
$ws.Cells.Item(85, 2).NumberFormat = "@"
$ws.Cells.Item(85, 2).Value = '4523'
$ws.Cells.Item(85, 2).Style = "Normal"
$ws.Cells.Item(85, 3).Value = 404
$ws.Cells.Item(86, 2).Value = '1911/5'
$ws.Cells.Item(86, 3).Value = 404
$ws.Cells.Item(87, 2).Value = '1912/5'
$ws.Cells.Item(87, 3).Value = 404
$ws.Cells.Item(88, 2).Value = '1912/5'
$ws.Cells.Item(88, 3).Value = 404

# The refreshed dataset has 3 fewer rows; drop the trailing now-empty rows
$ws.Rows(91).Delete()
$ws.Rows(90).Delete()
$ws.Rows(89).Delete()
